$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.446.79'
$ws.Range("E2").Value = '  +2.35%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.066.16'
$ws.Range("E3").Value = '  +2.31%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.43'
$ws.Range("E5").Value = '  -0.48%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.619'
$ws.Range("E6").Value = '  +2.94%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.56'
$ws.Range("E7").Value = '  +8.45%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.384'
$ws.Range("E9").Value = '  +3.84%  '

$ws.Range("E10").Value = '  +2.41%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0763'
$ws.Range("E11").Value = '  +1.81%  '

$ws.Range("E12").Value = '  +0.58%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.366.83'
$ws.Range("E13").Value = '  +2.17%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.41'
$ws.Range("E14").Value = '  +1.06%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.17'
$ws.Range("E15").Value = '  +4.82%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.777'
$ws.Range("E16").Value = '  +1.60%  '

$ws.Range("E17").Value = '  +1.37%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.065.11'
$ws.Range("E18").Value = '  +2.74%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '37.404.84'
$ws.Range("E19").Value = '  +2.53%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.23'
$ws.Range("E20").Value = '  +14.62%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '69.21'
$ws.Range("E21").Value = '  +2.20%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0812'
$ws.Range("E22").Value = '  +1.86%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '226.29'
$ws.Range("E23").Value = '  +2.67%  '

$ws.Range("E24").Value = '  -0.06%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.40'
$ws.Range("E25").Value = '  +0.82%  '

$ws.Range("E26").Value = '  +0.69%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '165.88'
$ws.Range("E27").Value = '  +1.53%  '

$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.47'
$ws.Range("E28").Value = '  +7.55%  '

$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.93'
$ws.Range("E29").Value = '  +3.39%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.128'
$ws.Range("E30").Value = '  -0.51%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '19.15'
$ws.Range("E31").Value = '  +1.32%  '

$ws.Range("E32").Value = '  +1.15%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.49'
$ws.Range("E33").Value = '  +2.66%  '

$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.59'
$ws.Range("E34").Value = '  +5.98%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0620'
$ws.Range("E35").Value = '  +3.07%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.57'
$ws.Range("E36").Value = '  +7.88%  '

$ws.Range("E37").Value = '  -0.15%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.94'
$ws.Range("E38").Value = '  +3.95%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.77'
$ws.Range("E39").Value = '  +0.11%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.31'
$ws.Range("E40").Value = '  -0.12%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.70'
$ws.Range("E41").Value = '  +11.95%  '

$ws.Range("E42").Value = '  -0.50%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0949'
$ws.Range("E43").Value = '  +0.34%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '96.77'
$ws.Range("E44").Value = '  +7.43%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.464.25'
$ws.Range("E45").Value = '  +0.53%  '

$ws.Range("E46").Value = '  +4.18%  '

$ws.Range("E47").Value = '  +5.91%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.73'
$ws.Range("E48").Value = '  +2.17%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.03'
$ws.Range("E49").Value = '  +2.31%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.24'
$ws.Range("E50").Value = '  +4.97%  '

$ws.Range("E51").Value = '  +1.73%  '
